$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the DEPTO value for the second row (João) from "SP" to "RS"
$ws.Range("C2").Value = "RS"

$wb.Save()
